$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 10 ("No.") : counts of right / wrong / not-attempted / max
# ------------------------------------------------------------------
$ws.Cells.Item(10, 1).Style = "mtitleStyle"
$ws.Cells.Item(10, 2).Value = 19
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(10, 4).Value = 7
$ws.Cells.Item(10, 5).Value = 28

# ------------------------------------------------------------------
# Row 11 ("Marking") : points per right / wrong answer
# ------------------------------------------------------------------
$ws.Cells.Item(11, 1).Style = "mtitleStyle"
$ws.Cells.Item(11, 2).Value = 4
$ws.Cells.Item(11, 3).Value = -1

# ------------------------------------------------------------------
# Row 12 ("Total") : totals and final score
# ------------------------------------------------------------------
$ws.Cells.Item(12, 1).Style = "mtitleStyle"
$ws.Cells.Item(12, 2).Value = 76
$ws.Cells.Item(12, 3).Value = -2
$ws.Cells.Item(12, 5).Value = "74/112"

# ------------------------------------------------------------------
# Row 15: the third "Student Ans"/"Correct Ans" header pair (G/H) is dropped
# ------------------------------------------------------------------
$ws.Range("G15:H15").Clear()

# ------------------------------------------------------------------
# Rows 16-40: student answer / correct answer grid.
# Column A = student's answer (blank if not attempted), colored by
# correctness; column B = the correct answer (reference column).
# Rows 16-18 still carry a second question pair in D/E; everything
# else (including the whole third G/H question pair) goes away.
# ------------------------------------------------------------------
    # Row 16
    $ws.Cells.Item(16, 1).Value = "Option A"
    $ws.Cells.Item(16, 1).Style = "correctStyle"
    $ws.Cells.Item(16, 2).Value = "Option A"
    $ws.Cells.Item(16, 2).Style = "absoluteStyle"
    $ws.Cells.Item(16, 4).Value = "Option A"
    $ws.Cells.Item(16, 4).Style = "correctStyle"
    $ws.Cells.Item(16, 5).Value = "Option A"
    $ws.Cells.Item(16, 5).Style = "absoluteStyle"
    # Row 17
    $ws.Cells.Item(17, 1).Value = ""
    $ws.Cells.Item(17, 1).Style = "normalStyle"
    $ws.Cells.Item(17, 2).Value = "Option D"
    $ws.Cells.Item(17, 2).Style = "absoluteStyle"
    $ws.Cells.Item(17, 4).Value = "Option C"
    $ws.Cells.Item(17, 4).Style = "correctStyle"
    $ws.Cells.Item(17, 5).Value = "Option C"
    $ws.Cells.Item(17, 5).Style = "absoluteStyle"
    # Row 18
    $ws.Cells.Item(18, 1).Value = "Option B"
    $ws.Cells.Item(18, 1).Style = "correctStyle"
    $ws.Cells.Item(18, 2).Value = "Option B"
    $ws.Cells.Item(18, 2).Style = "absoluteStyle"
    $ws.Cells.Item(18, 4).Value = "Option D"
    $ws.Cells.Item(18, 4).Style = "correctStyle"
    $ws.Cells.Item(18, 5).Value = "Option D"
    $ws.Cells.Item(18, 5).Style = "absoluteStyle"
    # Row 19
    $ws.Cells.Item(19, 1).Value = "Option C"
    $ws.Cells.Item(19, 1).Style = "correctStyle"
    $ws.Cells.Item(19, 2).Value = "Option C"
    $ws.Cells.Item(19, 2).Style = "absoluteStyle"
    # Row 20
    $ws.Cells.Item(20, 1).Value = "Option B"
    $ws.Cells.Item(20, 1).Style = "correctStyle"
    $ws.Cells.Item(20, 2).Value = "Option B"
    $ws.Cells.Item(20, 2).Style = "absoluteStyle"
    # Row 21
    $ws.Cells.Item(21, 1).Value = "Option C"
    $ws.Cells.Item(21, 1).Style = "correctStyle"
    $ws.Cells.Item(21, 2).Value = "Option C"
    $ws.Cells.Item(21, 2).Style = "absoluteStyle"
    # Row 22
    $ws.Cells.Item(22, 1).Value = "Option D"
    $ws.Cells.Item(22, 1).Style = "correctStyle"
    $ws.Cells.Item(22, 2).Value = "Option D"
    $ws.Cells.Item(22, 2).Style = "absoluteStyle"
    # Row 23
    $ws.Cells.Item(23, 1).Value = "Option D"
    $ws.Cells.Item(23, 1).Style = "correctStyle"
    $ws.Cells.Item(23, 2).Value = "Option D"
    $ws.Cells.Item(23, 2).Style = "absoluteStyle"
    # Row 24
    $ws.Cells.Item(24, 1).Value = ""
    $ws.Cells.Item(24, 1).Style = "normalStyle"
    $ws.Cells.Item(24, 2).Value = "Option A"
    $ws.Cells.Item(24, 2).Style = "absoluteStyle"
    # Row 25
    $ws.Cells.Item(25, 1).Value = "Option A"
    $ws.Cells.Item(25, 1).Style = "correctStyle"
    $ws.Cells.Item(25, 2).Value = "Option A"
    $ws.Cells.Item(25, 2).Style = "absoluteStyle"
    # Row 26
    $ws.Cells.Item(26, 1).Value = "Option D"
    $ws.Cells.Item(26, 1).Style = "incorrectStyle"
    $ws.Cells.Item(26, 2).Value = "Option C"
    $ws.Cells.Item(26, 2).Style = "absoluteStyle"
    # Row 27
    $ws.Cells.Item(27, 1).Value = "Option A"
    $ws.Cells.Item(27, 1).Style = "correctStyle"
    $ws.Cells.Item(27, 2).Value = "Option A"
    $ws.Cells.Item(27, 2).Style = "absoluteStyle"
    # Row 28
    $ws.Cells.Item(28, 1).Value = "Option D"
    $ws.Cells.Item(28, 1).Style = "correctStyle"
    $ws.Cells.Item(28, 2).Value = "Option D"
    $ws.Cells.Item(28, 2).Style = "absoluteStyle"
    # Row 29
    $ws.Cells.Item(29, 1).Value = ""
    $ws.Cells.Item(29, 1).Style = "normalStyle"
    $ws.Cells.Item(29, 2).Value = "Option D"
    $ws.Cells.Item(29, 2).Style = "absoluteStyle"
    # Row 30
    $ws.Cells.Item(30, 1).Value = "Option B"
    $ws.Cells.Item(30, 1).Style = "correctStyle"
    $ws.Cells.Item(30, 2).Value = "Option B"
    $ws.Cells.Item(30, 2).Style = "absoluteStyle"
    # Row 31
    $ws.Cells.Item(31, 1).Value = ""
    $ws.Cells.Item(31, 1).Style = "normalStyle"
    $ws.Cells.Item(31, 2).Value = "Option D"
    $ws.Cells.Item(31, 2).Style = "absoluteStyle"
    # Row 32
    $ws.Cells.Item(32, 1).Value = "Option C"
    $ws.Cells.Item(32, 1).Style = "correctStyle"
    $ws.Cells.Item(32, 2).Value = "Option C"
    $ws.Cells.Item(32, 2).Style = "absoluteStyle"
    # Row 33
    $ws.Cells.Item(33, 1).Value = "Option D"
    $ws.Cells.Item(33, 1).Style = "correctStyle"
    $ws.Cells.Item(33, 2).Value = "Option D"
    $ws.Cells.Item(33, 2).Style = "absoluteStyle"
    # Row 34
    $ws.Cells.Item(34, 1).Value = "Option B"
    $ws.Cells.Item(34, 1).Style = "correctStyle"
    $ws.Cells.Item(34, 2).Value = "Option B"
    $ws.Cells.Item(34, 2).Style = "absoluteStyle"
    # Row 35
    $ws.Cells.Item(35, 1).Value = ""
    $ws.Cells.Item(35, 1).Style = "normalStyle"
    $ws.Cells.Item(35, 2).Value = "Option D"
    $ws.Cells.Item(35, 2).Style = "absoluteStyle"
    # Row 36
    $ws.Cells.Item(36, 1).Value = "Option D"
    $ws.Cells.Item(36, 1).Style = "incorrectStyle"
    $ws.Cells.Item(36, 2).Value = "Option A"
    $ws.Cells.Item(36, 2).Style = "absoluteStyle"
    # Row 37
    $ws.Cells.Item(37, 1).Value = ""
    $ws.Cells.Item(37, 1).Style = "normalStyle"
    $ws.Cells.Item(37, 2).Value = "Option A"
    $ws.Cells.Item(37, 2).Style = "absoluteStyle"
    # Row 38
    $ws.Cells.Item(38, 1).Value = "Option A"
    $ws.Cells.Item(38, 1).Style = "correctStyle"
    $ws.Cells.Item(38, 2).Value = "Option A"
    $ws.Cells.Item(38, 2).Style = "absoluteStyle"
    # Row 39
    $ws.Cells.Item(39, 1).Value = "Option D"
    $ws.Cells.Item(39, 1).Style = "correctStyle"
    $ws.Cells.Item(39, 2).Value = "Option D"
    $ws.Cells.Item(39, 2).Style = "absoluteStyle"
    # Row 40
    $ws.Cells.Item(40, 1).Value = ""
    $ws.Cells.Item(40, 1).Style = "normalStyle"
    $ws.Cells.Item(40, 2).Value = "Option D"
    $ws.Cells.Item(40, 2).Style = "absoluteStyle"

# ------------------------------------------------------------------
# Drop everything that no longer exists in the new layout:
#  - the second question pair (D/E) for every row below 18
#  - the whole third question pair (G/H) for every data row
# ------------------------------------------------------------------
$ws.Range("D19:E40").Clear()
$ws.Range("G16:H40").Clear()
